# Generate Report for Handoff
# Adds two new rows (4c3df179-... and af41e3fc-...) to the Overview sheet
# and to the zh-cn / de-de detail sheets, mirroring the existing rows for
# 4cb39c9b-... / 57cda94f-....

$wb = $excel.ActiveWorkbook

$uuid1 = "4c3df179-2173-491f-9fe3-09b7473753bd"
$uuid2 = "af41e3fc-efcf-45dd-9817-8262d2f016a0"
$hash1 = "31cf03acfa0db16f5d278abccd9117a178cae439"
$hash2 = "20c5f16d8fd6e0e7af74920fc2d5bdd878697b13"

$md1 = "$uuid1.md"
$md2 = "$uuid2.md"
$xlf1zh = "$uuid1.$hash1.zh-cn.xlf"
$xlf2zh = "$uuid2.$hash2.zh-cn.xlf"
$xlf1de = "$uuid1.$hash1.de-de.xlf"
$xlf2de = "$uuid2.$hash2.de-de.xlf"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/6b4c8d3c7e0399cc9323512d38d64852f93d4545/e2e/$md1"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/6b4c8d3c7e0399cc9323512d38d64852f93d4545/e2e/$md2"

$xlfUrl1zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52f182f3e12731af020f962d39f313a607d5a0e3/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$xlf1zh"
$xlfUrl2zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52f182f3e12731af020f962d39f313a607d5a0e3/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$xlf2zh"

$xlfUrl1de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e95623598978c34c33e58f470ae4abadd50c263/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$xlf1de"
$xlfUrl2de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e95623598978c34c33e58f470ae4abadd50c263/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$xlf2de"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $md1
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-13-17 06:13:03"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl1, "", "", $md1)

$wsOverview.Range("A5").Value = $md2
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-13-17 06:13:03"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $mdUrl2, "", "", $md2)

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest
# Handoff File | Latest Handoff Datetime | ... | Handoff Reason |
# Dependency From
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = $md1
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = $xlf1zh
$wsZhCn.Range("E4").Value = "2016-03-17 06:12:56"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I4").Value = "Include"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $mdUrl1, "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B4"), $mdUrl1, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), $xlfUrl1zh, "", "", $xlf1zh)

$wsZhCn.Range("A5").Value = $md2
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = $xlf2zh
$wsZhCn.Range("E5").Value = "2016-03-17 06:12:56"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I5").Value = "Include"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), $mdUrl2, "", "", $md2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B5"), $mdUrl2, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), $xlfUrl2zh, "", "", $xlf2zh)

# ---------------------------------------------------------------------
# Sheet "de-de": same shape as "zh-cn"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = $md1
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = $xlf1de
$wsDeDe.Range("E4").Value = "2016-03-17 06:13:03"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I4").Value = "Include"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $mdUrl1, "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B4"), $mdUrl1, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), $xlfUrl1de, "", "", $xlf1de)

$wsDeDe.Range("A5").Value = $md2
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = $xlf2de
$wsDeDe.Range("E5").Value = "2016-03-17 06:13:03"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I5").Value = "Include"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), $mdUrl2, "", "", $md2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B5"), $mdUrl2, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), $xlfUrl2de, "", "", $xlf2de)
